# Applies cached Universalis market-data refresh values to the per-job
# "Leve Profits" tables (columns H:N) across all 8 crafting-job worksheets.
# Values were produced offline by the scheduled market-data runner and are
# written here as plain values (no formulas exist in these cells).
#
# Each data line is: SheetName,Row,Column,NewValue  (NewValue "CLEAR" means
# the cell must be removed/emptied rather than set).

$wb = $excel.ActiveWorkbook

$changes = @"
ALC,51,H,5531.4287
ALC,51,I,0
ALC,51,J,5531.4287
ALC,51,K,0
ALC,51,L,5531.4287
ALC,51,M,CLEAR
ALC,51,N,-6499.4287
ALC,74,H,8384.615
ALC,74,I,8187.375
ALC,74,J,8700.200000000001
ALC,74,K,8187.375
ALC,74,L,8700.200000000001
ALC,74,M,-7251.375
ALC,74,N,-10572.2
ALC,76,H,8750
ALC,76,I,0
ALC,76,J,8750
ALC,76,K,0
ALC,76,L,8750
ALC,76,M,CLEAR
ALC,76,N,-9380
ALC,77,H,8384.615
ALC,77,I,8187.375
ALC,77,J,8700.200000000001
ALC,77,K,40936.875
ALC,77,L,43501
ALC,77,M,-36256.875
ALC,77,N,-52861
ALC,79,H,8750
ALC,79,I,0
ALC,79,J,8750
ALC,79,K,0
ALC,79,L,8750
ALC,79,M,CLEAR
ALC,79,N,-10934
ALC,100,H,3399.077
ALC,100,I,3220
ALC,100,J,3802
ALC,100,K,3220
ALC,100,L,3802
ALC,100,M,-2679
ALC,100,N,-4884
ALC,124,H,99999
ALC,124,J,99999
ALC,124,L,99999
ALC,124,N,-109819
ALC,132,H,3012.4
ALC,132,I,1427.95
ALC,132,J,15688
ALC,132,K,4283.85
ALC,132,L,47064
ALC,132,M,-1753.85
ALC,132,N,-52124
ARM,2,H,2471.8
ARM,2,I,846.4545000000001
ARM,2,J,4458.3335
ARM,2,K,846.4545000000001
ARM,2,L,4458.3335
ARM,2,M,-733.4545000000001
ARM,2,N,-4684.3335
ARM,3,H,2999
ARM,3,I,2999
ARM,3,J,0
ARM,3,K,2999
ARM,3,L,0
ARM,3,M,-2884
ARM,3,N,CLEAR
ARM,61,H,4162.4614
ARM,61,I,2599.6667
ARM,61,J,5502
ARM,61,K,2599.6667
ARM,61,L,5502
ARM,61,M,-2387.6667
ARM,61,N,-5926
ARM,74,H,4932.8335
ARM,74,I,1899.25
ARM,74,K,1899.25
ARM,74,M,-1025.25
ARM,77,H,4932.8335
ARM,77,I,1899.25
ARM,77,K,9496.25
ARM,77,M,-5128.25
ARM,102,H,5639.5386
ARM,102,I,5739.125
ARM,102,J,4444.5
ARM,102,K,5739.125
ARM,102,L,4444.5
ARM,102,M,-4117.125
ARM,102,N,-7688.5
ARM,110,H,1233.6471
ARM,110,I,1206.5454
ARM,110,J,1283.3334
ARM,110,K,1206.5454
ARM,110,L,1283.3334
ARM,110,M,838.4546
ARM,110,N,-5373.3334
ARM,116,H,2471.8
ARM,116,I,846.4545000000001
ARM,116,J,4458.3335
ARM,116,K,846.4545000000001
ARM,116,L,4458.3335
ARM,116,M,1447.5455
ARM,116,N,-9046.333500000001
ARM,136,H,4162.4614
ARM,136,I,2599.6667
ARM,136,J,5502
ARM,136,K,7799.000100000001
ARM,136,L,16506
ARM,136,M,-5249.000100000001
ARM,136,N,-21606
BSM,3,H,2471.8
BSM,3,I,846.4545000000001
BSM,3,J,4458.3335
BSM,3,K,846.4545000000001
BSM,3,L,4458.3335
BSM,3,M,-732.4545000000001
BSM,3,N,-4686.3335
BSM,14,H,9998
BSM,14,I,9998
BSM,14,J,0
BSM,14,K,9998
BSM,14,L,0
BSM,14,M,-9826
BSM,14,N,CLEAR
BSM,38,H,30803.834
BSM,38,J,30500
BSM,38,L,30500
BSM,38,N,-31332
BSM,99,H,3393.7
BSM,99,I,3991.2856
BSM,99,K,3991.2856
BSM,99,M,-2493.2856
BSM,105,H,4825.75
BSM,105,I,4825.75
BSM,105,K,4825.75
BSM,105,M,-3078.75
BSM,107,H,2070.2856
BSM,107,I,1271.3636
BSM,107,K,1271.3636
BSM,107,M,648.6364000000001
CRP,16,H,1260.7142
CRP,16,I,1000
CRP,16,J,1304.1666
CRP,16,K,1000
CRP,16,L,1304.1666
CRP,16,M,-713
CRP,16,N,-1878.1666
CRP,31,H,2085.4038
CRP,31,I,786.1905
CRP,31,K,786.1905
CRP,31,M,-491.1905
CRP,34,H,2085.4038
CRP,34,I,786.1905
CRP,34,K,786.1905
CRP,34,M,-584.1905
CRP,113,H,1260.7142
CRP,113,I,1000
CRP,113,J,1304.1666
CRP,113,K,1000
CRP,113,L,1304.1666
CRP,113,M,1170
CRP,113,N,-5644.1666
CRP,114,H,80000
CRP,114,J,80000
CRP,114,L,80000
CRP,114,N,-88678
CUL,56,H,9777.700000000001
CUL,56,I,9777.700000000001
CUL,56,K,9777.700000000001
CUL,56,M,-9247.700000000001
CUL,63,H,15246.333
CUL,63,I,13425.5
CUL,63,K,40276.5
CUL,63,M,-39527.5
CUL,66,H,15246.333
CUL,66,I,13425.5
CUL,66,K,120829.5
CUL,66,M,-117085.5
CUL,82,H,19999
CUL,82,J,19999
CUL,82,L,59997
CUL,82,N,-60809
CUL,85,H,19999
CUL,85,J,19999
CUL,85,L,59997
CUL,85,N,-62805
CUL,92,H,202.875
CUL,92,I,168.75
CUL,92,J,237
CUL,92,K,506.25
CUL,92,L,711
CUL,92,M,741.75
CUL,92,N,-3207
CUL,107,H,626.1515000000001
CUL,107,I,667.8461
CUL,107,J,599.05
CUL,107,K,2003.5383
CUL,107,L,1797.15
CUL,107,M,-83.53829999999994
CUL,107,N,-5637.15
CUL,140,H,1252.72
CUL,140,I,1031.8572
CUL,140,J,2412.25
CUL,140,K,3095.5716
CUL,140,L,7236.75
CUL,140,M,2084.4284
CUL,140,N,-17596.75
GSM,70,H,43964.844
GSM,70,I,50381.285
GSM,70,J,25998.8
GSM,70,K,50381.285
GSM,70,L,25998.8
GSM,70,M,-50111.285
GSM,70,N,-26538.8
GSM,73,H,43964.844
GSM,73,I,50381.285
GSM,73,J,25998.8
GSM,73,K,50381.285
GSM,73,L,25998.8
GSM,73,M,-49445.285
GSM,73,N,-27870.8
GSM,107,H,1675.6154
GSM,107,I,995.5
GSM,107,J,2763.8
GSM,107,K,995.5
GSM,107,L,2763.8
GSM,107,M,924.5
GSM,107,N,-6603.8
LTW,68,H,5325.7144
LTW,68,I,1299.2
LTW,68,J,15392
LTW,68,K,1299.2
LTW,68,L,15392
LTW,68,M,-550.2
LTW,68,N,-16890
LTW,71,H,5325.7144
LTW,71,I,1299.2
LTW,71,J,15392
LTW,71,K,6496
LTW,71,L,76960
LTW,71,M,-2752
LTW,71,N,-84448
LTW,82,H,1142.3636
LTW,82,I,1161.875
LTW,82,J,1090.3334
LTW,82,K,1161.875
LTW,82,L,1090.3334
LTW,82,M,-800.875
LTW,82,N,-1812.3334
LTW,85,H,1142.3636
LTW,85,I,1161.875
LTW,85,J,1090.3334
LTW,85,K,1161.875
LTW,85,L,1090.3334
LTW,85,M,86.125
LTW,85,N,-3586.3334
LTW,100,H,2563.05
LTW,100,I,2633.2
LTW,100,J,2492.9
LTW,100,K,2633.2
LTW,100,L,2492.9
LTW,100,M,-2092.2
LTW,100,N,-3574.9
LTW,122,H,4932.8667
LTW,122,I,4731.7
LTW,122,J,5335.2
LTW,122,K,14195.1
LTW,122,L,16005.6
LTW,122,M,-11745.1
LTW,122,N,-20905.6
WVR,107,H,943.4194
WVR,107,I,679.5909
WVR,107,J,1588.3334
WVR,107,K,2038.7727
WVR,107,L,4765.0002
WVR,107,M,-118.7727
WVR,107,N,-8605.0002
WVR,113,H,1865.125
WVR,113,I,1887
WVR,113,J,1799.5
WVR,113,K,5661
WVR,113,L,5398.5
WVR,113,M,-3491
WVR,113,N,-9738.5
"@

$changeLines = $changes -split "`r?`n" | Where-Object { $_.Trim() -ne "" }

foreach ($line in $changeLines) {
    $parts = $line.Trim() -split ","
    $sheetName = $parts[0]
    $row       = $parts[1]
    $col       = $parts[2]
    $newValue  = $parts[3]

    $ws = $wb.Worksheets.Item($sheetName)
    $cell = $ws.Range("$col$row")

    if ($newValue -eq "CLEAR") {
        $cell.ClearContents()
    } else {
        $cell.Value = [double]$newValue
    }
}
